$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9227734804153442
$ws.Range("B1").Value = 1.101991057395935
$ws.Range("C1").Value = 0.9133025407791138
$ws.Range("D1").Value = 0.8500090837478638
$ws.Range("E1").Value = 0.8839741945266724
